$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '''26.564.83'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -3.07%  '
$c = $ws.Range("D3")
$c.Value = '''1.806.59'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -2.80%  '
$c = $ws.Range("D4")
$c.Value = '''1.009'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.50%  '
$c = $ws.Range("D5")
$c.Value = '''1.009'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '
$c = $ws.Range("D6")
$c.Value = '''308.78'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.91%  '
$c = $ws.Range("D7")
$c.Value = '''0.4545'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.92%  '
$c = $ws.Range("D8")
$c.Value = '''0.3666'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.61%  '
$c = $ws.Range("D9")
$c.Value = '''0.07137'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.56%  '
$c = $ws.Range("D10")
$c.Value = '''0.8719'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("E11").Value = '  -0.65%  '
$c = $ws.Range("D12")
$c.Value = '''19.24'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.77%  '
$c = $ws.Range("D13")
$c.Value = '''1.847.29'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.30%  '
$c = $ws.Range("D14")
$c.Value = '''5.287'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.16%  '
$c = $ws.Range("D15")
$c.Value = '''6.332'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.67%  '
$c = $ws.Range("D16")
$c.Value = '''86.66'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -5.74%  '
$c = $ws.Range("D17")
$c.Value = '''1.009'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.48%  '
$c = $ws.Range("D18")
$c.Value = '''0.000008581'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -4.52%  '
$ws.Range("E19").Value = '  +0.54%  '
$c = $ws.Range("D20")
$c.Value = '''26.602.30'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.01%  '
$c = $ws.Range("D21")
$c.Value = '''14.26'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.64%  '
$c = $ws.Range("D22")
$c.Value = '''4.961'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.38%  '
$c = $ws.Range("D23")
$c.Value = '''2.049.04'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.35%  '
$c = $ws.Range("D24")
$c.Value = '''10.36'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.86%  '
$c = $ws.Range("D25")
$c.Value = '''1.982'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.16%  '
$c = $ws.Range("D26")
$c.Value = '''151.08'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.83%  '
$c = $ws.Range("D27")
$c.Value = '''17.92'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.95%  '
$c = $ws.Range("D28")
$c.Value = '''1.998'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.06%  '
$c = $ws.Range("D29")
$c.Value = '''113.03'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.96%  '
$c = $ws.Range("D30")
$c.Value = '''4.878'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -4.43%  '
$c = $ws.Range("D31")
$c.Value = '''0.08697'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.79%  '
$c = $ws.Range("D32")
$c.Value = '''3.067'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.09%  '
$c = $ws.Range("D33")
$c.Value = '''0.7346'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.64%  '
$c = $ws.Range("D34")
$c.Value = '''4.439'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.74%  '
$c = $ws.Range("D35")
$c.Value = '''1.115'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -5.34%  '
$c = $ws.Range("D36")
$c.Value = '''2.505'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -7.41%  '
$c = $ws.Range("D37")
$c.Value = '''1.078'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.38%  '
$c = $ws.Range("D38")
$c.Value = '''0.01920'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.26%  '
$c = $ws.Range("D39")
$c.Value = '''0.05098'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.80%  '
$c = $ws.Range("D40")
$c.Value = '''2.864'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.22%  '
$c = $ws.Range("D41")
$c.Value = '''6.879'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.28%  '
$c = $ws.Range("D42")
$c.Value = '''0.4916'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -4.67%  '
$c = $ws.Range("D43")
$c.Value = '''0.1570'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -4.38%  '
$c = $ws.Range("D44")
$c.Value = '''8.133'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.18%  '
$c = $ws.Range("D45")
$c.Value = '''1.010'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.56%  '
$c = $ws.Range("D46")
$c.Value = '''0.4600'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.58%  '
$c = $ws.Range("D47")
$c.Value = '''102.17'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.25%  '
$c = $ws.Range("D48")
$c.Value = '''9.938'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.98%  '
$c = $ws.Range("D49")
$c.Value = '''1.584'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -4.23%  '
$c = $ws.Range("D50")
$c.Value = '''0.05998'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -3.61%  '
$c = $ws.Range("D51")
$c.Value = '''63.63'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.03%  '
